$wb = $excel.ActiveWorkbook

# --- Sheet "Hazard-focused" (sheet1) ---
$ws1 = $wb.Worksheets.Item("Hazard-focused")

# Header row 1
$ws1.Range("A1").Value = "Hazard Noun/Subject"
$ws1.Range("B1").Value = "Action/Descriptor"
$ws1.Range("B1").Font.Bold = $true

# Row 2: split "resource, crew; limited, share, lack, fatigue"
$ws1.Range("A2").Value = "resource, crew"
$ws1.Range("B2").Value = "limited, share, lack, fatigue"

# Row 3: split "highway, road; close, closure"
$ws1.Range("A3").Value = "highway, road"
$ws1.Range("B3").Value = "close, closure"

# Row 4: split "ground; aircraft, helicopter, heli, copter, aerial, tanker"
$ws1.Range("A4").Value = "ground"
$ws1.Range("B4").Value = "aircraft, helicopter, heli, copter, aerial, tanker"

$ws1.Range("J13").Select()
